$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.26
$ws.Range("J3").Value = 3.25
$ws.Range("S3").Value = 2.62
$ws.Range("U3").Value = 2.24
$ws.Range("I4").Value = 1.24
$ws.Range("N4").Value = 7.8
$ws.Range("O4").Value = 1.11
$ws.Range("R4").Value = 1.95
$ws.Range("S4").Value = 1.77
$ws.Range("T4").Value = 1.92
$ws.Range("U4").Value = 1.89
$ws.Range("V4").Value = 5.1
$ws.Range("X4").Value = 60
$ws.Range("Y4").Value = 17
$ws.Range("Z4").Value = 12.5
$ws.Range("AA4").Value = 12.5
$ws.Range("AB4").Value = 75
$ws.Range("AC4").Value = 26
$ws.Range("AD4").Value = 15.5
$ws.Range("AE4").Value = 16.5
$ws.Range("AF4").Value = 190
$ws.Range("AG4").Value = 65
$ws.Range("AH4").Value = 44
$ws.Range("AI4").Value = 42
$ws.Range("AK4").Value = 260
$ws.Range("AL4").Value = 170
$ws.Range("AM4").Value = 180
$ws.Range("AN4").Value = 260
$ws.Range("AO4").Value = 3.5
$ws.Range("R5").Value = 1.63
$ws.Range("U5").Value = 2.68
$ws.Range("F6").Value = 3.25
$ws.Range("G6").Value = 3.35
$ws.Range("I6").Value = 2.22
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 1.54
$ws.Range("U6").Value = 2.82
$ws.Range("AD6").Value = 11.5
$ws.Range("AH6").Value = 14.5
$ws.Range("F8").Value = 3.9
$ws.Range("G8").Value = 4.2
$ws.Range("H8").Value = 1.93
$ws.Range("I8").Value = 2.04
$ws.Range("J8").Value = 3.8
$ws.Range("K8").Value = 4.2
$ws.Range("I9").Value = 5.4
$ws.Range("J9").Value = 3.5
$ws.Range("F10").Value = 2.34
$ws.Range("G10").Value = 2.8
$ws.Range("H10").Value = 2.86
$ws.Range("J10").Value = 3.2
$ws.Range("P10").Value = 2.08
$ws.Range("Q10").Value = 1.73
$ws.Range("H11").Value = 1.09
$ws.Range("J11").Value = 4.1
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 3.85
$ws.Range("N12").Value = 3.75
$ws.Range("P12").Value = 1.94
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.36
$ws.Range("U12").Value = 2.04
$ws.Range("X12").Value = 14
$ws.Range("Y12").Value = 17
$ws.Range("AC12").Value = 8.4
$ws.Range("AD12").Value = 20
$ws.Range("F13").Value = 9.6
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 1.33
$ws.Range("I13").Value = 1.34
$ws.Range("Q13").Value = 1.38
$ws.Range("R13").Value = 1.99
$ws.Range("S13").Value = 1.95
$ws.Range("T13").Value = 1.69
$ws.Range("U13").Value = 2.36
$ws.Range("X13").Value = 46
$ws.Range("AA13").Value = 13
$ws.Range("AL13").Value = 90
$ws.Range("G14").Value = 3.9
$ws.Range("H14").Value = 2.04
$ws.Range("I14").Value = 2.06
$ws.Range("P14").Value = 2.38
$ws.Range("Q14").Value = 1.69
$ws.Range("R14").Value = 1.55
$ws.Range("S14").Value = 2.7
$ws.Range("U14").Value = 2.48
$ws.Range("AC14").Value = 9.199999999999999
$ws.Range("AI14").Value = 28
$ws.Range("AJ14").Value = 70
$ws.Range("N15").Value = 3.6
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 1.88
$ws.Range("Q15").Value = 2.1
$ws.Range("S15").Value = 3.75
$ws.Range("T15").Value = 1.93
$ws.Range("U15").Value = 2.02
$ws.Range("X15").Value = 13.5
$ws.Range("Y15").Value = 16
$ws.Range("Z15").Value = 38
$ws.Range("AB15").Value = 8.4
$ws.Range("F16").Value = 1.48
$ws.Range("G16").Value = 1.5
$ws.Range("H16").Value = 6.6
$ws.Range("I16").Value = 7
$ws.Range("R16").Value = 1.91
$ws.Range("S16").Value = 2.04
$ws.Range("AC16").Value = 14
$ws.Range("AD16").Value = 27
$ws.Range("AK16").Value = 13.5
$ws.Range("H17").Value = 24
$ws.Range("I17").Value = 25
$ws.Range("U17").Value = 1.79
$ws.Range("Y17").Value = 90
$ws.Range("AD17").Value = 90
$ws.Range("AG17").Value = 14
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 28
$ws.Range("K18").Value = 13
$ws.Range("U18").Value = 1.88
$ws.Range("X18").Value = 75
$ws.Range("AC18").Value = 32
$ws.Range("AG18").Value = 15.5
$ws.Range("AH18").Value = 48
$ws.Range("AJ18").Value = 9.800000000000001
$ws.Range("AL18").Value = 100
